# Auto-generated edit script
# Applies updated market-price figures (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) for a batch of leve rows across several crafting-job sheets, as
# pulled from the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 1685.0125
$ws.Range("J17").Value = 1835.2253
$ws.Range("L17").Value = 5505.6759
$ws.Range("N17").Value = -5841.6759
# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 23544.455
$ws.Range("I88").Value = 7664.3335
$ws.Range("J88").Value = 29499.5
$ws.Range("K88").Value = 7664.3335
$ws.Range("L88").Value = 29499.5
$ws.Range("M88").Value = -7258.3335
$ws.Range("N88").Value = -30311.5
# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 23544.455
$ws.Range("I91").Value = 7664.3335
$ws.Range("J91").Value = 29499.5
$ws.Range("K91").Value = 7664.3335
$ws.Range("L91").Value = 29499.5
$ws.Range("M91").Value = -6260.3335
$ws.Range("N91").Value = -32307.5
# Row 97: Materia Worth / Potent Spiritbond Potion
$ws.Range("H97").Value = 1563.25
$ws.Range("J97").Value = 1518
$ws.Range("L97").Value = 4554
$ws.Range("N97").Value = -5546

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 285.41666
$ws.Range("I5").Value = 214.22223
$ws.Range("J5").Value = 499
$ws.Range("K5").Value = 214.22223
$ws.Range("L5").Value = 499
$ws.Range("M5").Value = -102.22223
$ws.Range("N5").Value = -723
# Row 46: Get Me the Usual / Heavy Steel Flanchard
$ws.Range("H46").Value = 28927.285
$ws.Range("J46").Value = 16698.2
$ws.Range("L46").Value = 16698.2
$ws.Range("N46").Value = -17336.2
# Row 88: The Mast Chance / Adamantite Rivets
$ws.Range("H88").Value = 3998.3635
$ws.Range("J88").Value = 3247
$ws.Range("L88").Value = 3247
$ws.Range("N88").Value = -4059
# Row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws.Range("H91").Value = 3998.3635
$ws.Range("J91").Value = 3247
$ws.Range("L91").Value = 3247
$ws.Range("N91").Value = -6055
# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 297.92307
$ws.Range("I97").Value = 314.63635
$ws.Range("J97").Value = 206
$ws.Range("K97").Value = 314.63635
$ws.Range("L97").Value = 206
$ws.Range("M97").Value = 181.36365
$ws.Range("N97").Value = -1198
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 1481.8636
$ws.Range("I122").Value = 812.86664
$ws.Range("K122").Value = 2438.59992
$ws.Range("M122").Value = 11.40008000000034
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 3843.8845
$ws.Range("I132").Value = 3095.4736
$ws.Range("K132").Value = 9286.4208
$ws.Range("M132").Value = -6756.4208

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 285.41666
$ws.Range("I4").Value = 214.22223
$ws.Range("J4").Value = 499
$ws.Range("K4").Value = 214.22223
$ws.Range("L4").Value = 499
$ws.Range("M4").Value = -99.22223
$ws.Range("N4").Value = -729
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 48145.57
$ws.Range("I105").Value = 82006.664
$ws.Range("J105").Value = 22749.75
$ws.Range("K105").Value = 82006.664
$ws.Range("L105").Value = 22749.75
$ws.Range("M105").Value = -80259.664
$ws.Range("N105").Value = -26243.75
# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 847.3158
$ws.Range("I107").Value = 364.64706
$ws.Range("K107").Value = 364.64706
$ws.Range("M107").Value = 1555.35294
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2761.0952
$ws.Range("I134").Value = 1735.3158
$ws.Range("K134").Value = 5205.9474
$ws.Range("M134").Value = -2670.9474

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7: It's Always Sunny in Vylbrand / Raisins
$ws.Range("H7").Value = 50182.1
$ws.Range("I7").Value = 177.375
$ws.Range("K7").Value = 532.125
$ws.Range("M7").Value = -420.125
# Row 17: Chew the Fat / Grilled Dodo
$ws.Range("H17").Value = 631.2857
$ws.Range("J17").Value = 100
$ws.Range("L17").Value = 300
$ws.Range("N17").Value = -638
# Row 19: The Bango Zango Diet / Parsnip Salad
$ws.Range("H19").Value = 4
$ws.Range("J19").Value = 4
$ws.Range("L19").Value = 12
$ws.Range("N19").Value = -360
# Row 23: Sweet Smell of Success / Lavender Oil
$ws.Range("H23").Value = 628
$ws.Range("I23").Value = 200
$ws.Range("J23").Value = 689.1429000000001
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 2067.4287
$ws.Range("M23").Value = -365
$ws.Range("N23").Value = -2537.4287
# Row 25: Flakes for Friends / Apple Tart
$ws.Range("H25").Value = 2124.0833
$ws.Range("I25").Value = 1563
$ws.Range("K25").Value = 4689
$ws.Range("M25").Value = -4520
# Row 30: Picnic Panic / Apple Tart
$ws.Range("H30").Value = 2124.0833
$ws.Range("I30").Value = 1563
$ws.Range("K30").Value = 4689
$ws.Range("M30").Value = -4587
# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 4497.143
$ws.Range("J68").Value = 4497.143
$ws.Range("L68").Value = 13491.429
$ws.Range("N68").Value = -15113.429
# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 4497.143
$ws.Range("J71").Value = 4497.143
$ws.Range("L71").Value = 40474.287
$ws.Range("N71").Value = -48586.287
# Row 103: West Meats East / Nomad Meat Pie
$ws.Range("H103").Value = 1301.7273
$ws.Range("I103").Value = 845
$ws.Range("J103").Value = 1849.8
$ws.Range("K103").Value = 2535
$ws.Range("L103").Value = 5549.4
$ws.Range("M103").Value = -1656
$ws.Range("N103").Value = -7307.4
# Row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Range("H121").Value = 3681.45
$ws.Range("I121").Value = 1846
$ws.Range("J121").Value = 4293.2666
$ws.Range("K121").Value = 5538
$ws.Range("L121").Value = 12879.7998
$ws.Range("M121").Value = -4228
$ws.Range("N121").Value = -15499.7998
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 23150394
$ws.Range("J131").Value = 20834608
$ws.Range("L131").Value = 62503824
$ws.Range("N131").Value = -62513904
# Row 138: Bring Me Your Tacos / Tacos Al Pastor
$ws.Range("H138").Value = 4465.2856
$ws.Range("I138").Value = 1042.8334
$ws.Range("K138").Value = 3128.5002
$ws.Range("M138").Value = 2011.4998

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 34: All Booked Up / Silver Magnifiers
$ws.Range("H34").Value = 40000
$ws.Range("I34").Value = 40000
$ws.Range("K34").Value = 40000
$ws.Range("M34").Value = -39732
# Row 76: The Monuments Mages / Hardsilver Magnifiers of Casting
$ws.Range("H76").Value = 40000
$ws.Range("I76").Value = 40000
$ws.Range("K76").Value = 40000
$ws.Range("M76").Value = -39685
# Row 79: Deal with It (L) / Hardsilver Magnifiers of Casting
$ws.Range("H79").Value = 40000
$ws.Range("I79").Value = 40000
$ws.Range("K79").Value = 40000
$ws.Range("M79").Value = -38908
# Row 123: Workplace Workout / Ametrine Ring of Fending
$ws.Range("H123").Value = 38266
$ws.Range("J123").Value = 38266
$ws.Range("L123").Value = 38266
$ws.Range("N123").Value = -43166

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 31: Open to Attack / Goatskin Jacket
$ws.Range("H31").Value = 12759744
$ws.Range("I31").Value = 523.5
$ws.Range("J31").Value = 25518964
$ws.Range("K31").Value = 523.5
$ws.Range("L31").Value = 25518964
$ws.Range("M31").Value = -275.5
$ws.Range("N31").Value = -25519460
# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 14727.444
$ws.Range("I40").Value = 14192.75
$ws.Range("J40").Value = 19005
$ws.Range("K40").Value = 14192.75
$ws.Range("L40").Value = 19005
$ws.Range("M40").Value = -14056.75
$ws.Range("N40").Value = -19277
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 5033.905
$ws.Range("I122").Value = 4682.467
$ws.Range("J122").Value = 5912.5
$ws.Range("K122").Value = 14047.401
$ws.Range("L122").Value = 17737.5
$ws.Range("M122").Value = -11597.401
$ws.Range("N122").Value = -22637.5
# Row 129: Loving Soles / Kumbhiraskin Boots of Gathering
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 1163.2
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1163.2
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3489.6
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -7829.6
# Row 129: Lifetime of Gleaning / Scarlet Moko Beret of Gathering
$ws.Range("H129").Value = 63000
$ws.Range("J129").Value = 63000
$ws.Range("L129").Value = 63000
$ws.Range("N129").Value = -73000
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 2394.8293
$ws.Range("I136").Value = 2012.3846
$ws.Range("J136").Value = 9852.5
$ws.Range("K136").Value = 6037.1538
$ws.Range("L136").Value = 29557.5
$ws.Range("M136").Value = -3487.1538
$ws.Range("N136").Value = -34657.5
